# Applies the "Add files via upload" commit to 3_6_2nd_order.xlsx:
#  - renames sheets "5_Quant" -> "5_" and "6_MultAns" -> "6_"
#  - replaces the quiz content of those two sheets with new questions
#    about terminal velocity / drag, and resizes/re-selects them
#  - moves the active tab from "4_" to the renamed "6_" sheet
#  - tidies the selection left behind on "4_"

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename sheets 5_Quant -> 5_ and 6_MultAns -> 6_
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("5_Quant")
$ws5.Name = "5_"

$ws6 = $wb.Worksheets.Item("6_MultAns")
$ws6.Name = "6_"

# ------------------------------------------------------------------
# 2. Rebuild "5_" (was 5_Quant) with the new quiz content
# ------------------------------------------------------------------
$ws5.Cells.Clear()

$ws5.Range("A1").Value = "What assumption allowed us to find the coefficient of drag for the falling penny?"
$ws5.Range("B1").Value = "Correct"
$ws5.Range("C1").Value = "Comment"

$ws5.Range("A2").Value = "We assumed a value for gravitational acceleration"
$ws5.Range("B2").Value = "N"
$ws5.Range("C2").Value = "This value we know pretty accurately."

$ws5.Range("A3").Value = "We assumed that the penny was not whirling while it fell"
$ws5.Range("B3").Value = "N"

$ws5.Range("A4").Value = "We assumed that the penny fell with its faces pointing up and down"
$ws5.Range("B4").Value = "N"
$ws5.Range("C4").Value = "We could make this assumption to find a C_d, but we didn't do that here"

$ws5.Range("A5").Value = "We assumed a value for terminal velocity"
$ws5.Range("B5").Value = "Y"
$ws5.Range("C5").Value = "Yep!  Assuming v_term allowed us to solve for C_d."

$ws5.Rows.Item(1).RowHeight = 60
$ws5.Rows.Item(2).RowHeight = 30
$ws5.Rows.Item(3).RowHeight = 45
$ws5.Rows.Item(4).RowHeight = 45
$ws5.Rows.Item(5).RowHeight = 30

$ws5.Columns.Item(1).ColumnWidth = 24.25

$ws5.Range("A1:C5").Select()

# ------------------------------------------------------------------
# 3. Rebuild "6_" (was 6_MultAns) with the new quiz content
# ------------------------------------------------------------------
$ws6.Cells.Clear()

$ws6.Range("A1").Value = "The position plot is mostly linear.  But in the previous notebook (without drag), the position plot was parabolic.  What causes the difference?"
$ws6.Range("B1").Value = "Correct"
$ws6.Range("C1").Value = "Comment"

$ws6.Range("A2").Value = "When we include drag, velocity is constant for much of the fall"
$ws6.Range("B2").Value = "Y"
$ws6.Range("C2").Value = "Yep!  There is no acceleration once the penny reaches terminal velocity, so the velocity is constant."

$ws6.Range("A3").Value = "When we include drag, acceleration is constant for much of the fall"
$ws6.Range("B3").Value = "N"
$ws6.Range("C3").Value = "Acceleration is constant when there is *no* drag (-9.8 m/s2)"

$ws6.Range("A4").Value = "When we include drag, position is constant for much of the fall"
$ws6.Range("B4").Value = "N"

$ws6.Rows.Item(1).RowHeight = 60
$ws6.Rows.Item(2).RowHeight = 45
$ws6.Rows.Item(3).RowHeight = 30
$ws6.Rows.Item(4).RowHeight = 30

$ws6.Range("C3").Select()
$ws6.Range("C3:C3").Select()

# ------------------------------------------------------------------
# 4. "6_" becomes the active tab (was "4_")
# ------------------------------------------------------------------
$ws6.Select()

# ------------------------------------------------------------------
# 5. Tidy up the selection left on "4_"
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("4_")
$ws4.Range("A1:C4").Select()

# Restore the intended active sheet/cell on "6_" (selecting "4_" above
# moved the active sheet away from it).
$ws6.Select()
$ws6.Range("C3").Select()
